# Update the "想去人数" (want-to-go count) figures in column F of both the
# "展览" and "全部类型" worksheets. These two sheets mirror each other's data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 81
    3  = 3098
    5  = 2639
    7  = 135
    9  = 1375
    11 = 60
    12 = 15
    13 = 1187
    14 = 357
    17 = 32
    19 = 73
    21 = 2516
    22 = 31
    23 = 283
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
